$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1893687707641196
$ws.Range("C2").Value = 0.5448504983388704
$ws.Range("J2").Value = 0.026578073089701
$ws.Range("P2").Value = 0.1262458471760797
$ws.Range("S2").Value = 0.1129568106312292
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.006134969325153374
$ws.Range("J3").Value = 0.03067484662576687
$ws.Range("P3").Value = 0.7423312883435583
$ws.Range("S3").Value = 0.2147239263803681
$ws.Range("J4").Value = 0.04878048780487805
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("B6").Value = 0.04700854700854701
$ws.Range("D6").Value = 0.0170940170940171
$ws.Range("F6").Value = 0.04700854700854701
$ws.Range("J6").Value = 0.2735042735042735
$ws.Range("O6").Value = 0.03846153846153846
$ws.Range("Q6").Value = 0.2136752136752137
$ws.Range("R6").Value = 0.05128205128205128
$ws.Range("S6").Value = 0.311965811965812
$ws.Range("B7").Value = 0.1339285714285714
$ws.Range("D7").Value = 0.01785714285714286
$ws.Range("F7").Value = 0.03125
$ws.Range("J7").Value = 0.1160714285714286
$ws.Range("O7").Value = 0.01785714285714286
$ws.Range("Q7").Value = 0.1741071428571428
$ws.Range("R7").Value = 0.08035714285714286
$ws.Range("S7").Value = 0.4285714285714285
$ws.Range("B8").Value = 0.08316831683168317
$ws.Range("D8").Value = 0.02376237623762376
$ws.Range("F8").Value = 0.04752475247524753
$ws.Range("J8").Value = 0.1386138613861386
$ws.Range("O8").Value = 0.03168316831683168
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.07524752475247524
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.07949790794979079
$ws.Range("D9").Value = 0.02510460251046025
$ws.Range("F9").Value = 0.06694560669456066
$ws.Range("J9").Value = 0.07949790794979079
$ws.Range("O9").Value = 0.01255230125523013
$ws.Range("Q9").Value = 0.2092050209205021
$ws.Range("R9").Value = 0.1171548117154812
$ws.Range("S9").Value = 0.4100418410041841
$ws.Range("B10").Value = 0.1108374384236453
$ws.Range("D10").Value = 0.01642036124794746
$ws.Range("E10").Value = 0.0008210180623973727
$ws.Range("F10").Value = 0.09277504105090312
$ws.Range("J10").Value = 0.1133004926108374
$ws.Range("O10").Value = 0.0180623973727422
$ws.Range("Q10").Value = 0.19376026272578
$ws.Range("R10").Value = 0.09359605911330049
$ws.Range("S10").Value = 0.3604269293924466
$ws.Range("G11").Value = 0.1331168831168831
$ws.Range("J11").Value = 0.08116883116883117
$ws.Range("K11").Value = 0.1623376623376623
$ws.Range("L11").Value = 0.6136363636363636
$ws.Range("S11").Value = 0.00974025974025974
$ws.Range("G12").Value = 0.7673267326732673
$ws.Range("J12").Value = 0.1633663366336634
$ws.Range("K12").Value = 0.01485148514851485
$ws.Range("L12").Value = 0.04455445544554455
$ws.Range("S12").Value = 0.009900990099009901
$ws.Range("G13").Value = 0.6851851851851852
$ws.Range("J13").Value = 0.2592592592592592
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.03846153846153846
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.06837606837606838
$ws.Range("J15").Value = 0.2991452991452991
$ws.Range("K15").Value = 0.06837606837606838
$ws.Range("M15").Value = 0.0170940170940171
$ws.Range("O15").Value = 0.04273504273504274
$ws.Range("S15").Value = 0.2991452991452991
$ws.Range("F16").Value = 0.01081081081081081
$ws.Range("H16").Value = 0.2054054054054054
$ws.Range("I16").Value = 0.1027027027027027
$ws.Range("J16").Value = 0.3243243243243243
$ws.Range("K16").Value = 0.1027027027027027
$ws.Range("M16").Value = 0.02162162162162162
$ws.Range("O16").Value = 0.06486486486486487
$ws.Range("S16").Value = 0.1675675675675676
$ws.Range("F17").Value = 0.02320675105485232
$ws.Range("H17").Value = 0.20042194092827
$ws.Range("I17").Value = 0.1223628691983122
$ws.Range("J17").Value = 0.3628691983122363
$ws.Range("K17").Value = 0.1139240506329114
$ws.Range("M17").Value = 0.02109704641350211
$ws.Range("O17").Value = 0.05274261603375528
$ws.Range("S17").Value = 0.1033755274261603
$ws.Range("F18").Value = 0.009523809523809525
$ws.Range("H18").Value = 0.1714285714285714
$ws.Range("I18").Value = 0.1095238095238095
$ws.Range("J18").Value = 0.3952380952380952
$ws.Range("K18").Value = 0.09523809523809523
$ws.Range("M18").Value = 0.004761904761904762
$ws.Range("O18").Value = 0.07142857142857142
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01278195488721805
$ws.Range("H19").Value = 0.2278195488721804
$ws.Range("I19").Value = 0.09398496240601503
$ws.Range("J19").Value = 0.3360902255639098
$ws.Range("K19").Value = 0.1052631578947368
$ws.Range("M19").Value = 0.02781954887218045
$ws.Range("N19").Value = 0.001503759398496241
$ws.Range("O19").Value = 0.06691729323308271
$ws.Range("S19").Value = 0.1278195488721804
